# Atualizacao de bases das ligas (Serbia Super Liga), do dia: 08-04-2024 as 21:28
# Applies row-level odds/result corrections + a HomeTeam/AwayTeam shared-string swap
# (Spartak Subotica <-> FK Cukaricki) as captured by the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 2
$ws.Range("B2").Value = 6983357
$ws.Range("G2").Value = "Radnicki Nis"
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 1.571
$ws.Range("L2").Value = 3.6
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = 1.533
$ws.Range("O2").Value = 3.6
$ws.Range("P2").Value = 5.25
$ws.Range("Q2").Value = -1
$ws.Range("R2").Value = 2
$ws.Range("S2").Value = 1.8
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 1.925
$ws.Range("V2").Value = 1.875
$ws.Range("W2").Value = 0.5329999999999999
$ws.Range("Z2").Value = 1
$ws.Range("AB2").Value = -1
$ws.Range("AC2").Value = 0.875

# Row 3
$ws.Range("B3").Value = 6979406
$ws.Range("G3").Value = "IMT Novi Belgrade"
$ws.Range("I3").Value = 1
$ws.Range("K3").Value = 1.909
$ws.Range("L3").Value = 3.2
$ws.Range("M3").Value = 3.6
$ws.Range("N3").Value = 1.909
$ws.Range("O3").Value = 3.1
$ws.Range("P3").Value = 3.8
$ws.Range("Q3").Value = -0.5
$ws.Range("R3").Value = 1.975
$ws.Range("S3").Value = 1.825
$ws.Range("T3").Value = 2.25
$ws.Range("U3").Value = 1.9
$ws.Range("V3").Value = 1.9
$ws.Range("W3").Value = 0.909
$ws.Range("Z3").Value = 0.9750000000000001
$ws.Range("AB3").Value = 0.8999999999999999
$ws.Range("AC3").Value = -1

# Row 14
$ws.Range("G14").Value = "Spartak Subotica"

# Row 17
$ws.Range("G17").Value = "FK Cukaricki"

# Row 21
$ws.Range("F21").Value = "Spartak Subotica"

# Row 24
$ws.Range("F24").Value = "FK Cukaricki"

# Row 26
$ws.Range("G26").Value = "Spartak Subotica"

# Row 28
$ws.Range("G28").Value = "FK Cukaricki"

# Row 31
$ws.Range("B31").Value = 6979427
$ws.Range("F31").Value = "Javor Ivanjica"
$ws.Range("G31").Value = "FK Backa Topola"
$ws.Range("I31").Value = 3
$ws.Range("K31").Value = 5.25
$ws.Range("L31").Value = 4
$ws.Range("M31").Value = 1.533
$ws.Range("N31").Value = 4
$ws.Range("O31").Value = 3.3
$ws.Range("P31").Value = 1.85
$ws.Range("Q31").Value = 0.5
$ws.Range("R31").Value = 1.875
$ws.Range("S31").Value = 1.925
$ws.Range("T31").Value = 2.25
$ws.Range("U31").Value = 1.825
$ws.Range("V31").Value = 1.975
$ws.Range("Y31").Value = 0.8500000000000001
$ws.Range("AA31").Value = 0.925
$ws.Range("AB31").Value = 0.825
$ws.Range("AC31").Value = -1

# Row 32
$ws.Range("B32").Value = 6979431
$ws.Range("F32").Value = "FK Napredak"
$ws.Range("G32").Value = "Mladost Lucani"
$ws.Range("I32").Value = 1
$ws.Range("K32").Value = 1.7
$ws.Range("L32").Value = 3.75
$ws.Range("M32").Value = 4.5
$ws.Range("N32").Value = 1.909
$ws.Range("O32").Value = 3.1
$ws.Range("P32").Value = 4
$ws.Range("Q32").Value = -0.5
$ws.Range("R32").Value = 1.975
$ws.Range("S32").Value = 1.825
$ws.Range("T32").Value = 2
$ws.Range("U32").Value = 1.875
$ws.Range("V32").Value = 1.925
$ws.Range("Y32").Value = 3
$ws.Range("AA32").Value = 0.825
$ws.Range("AB32").Value = -1
$ws.Range("AC32").Value = 0.925

# Row 38
$ws.Range("B38").Value = 6979433
$ws.Range("G38").Value = "FK Napredak"
$ws.Range("H38").Value = 1
$ws.Range("I38").Value = 3
$ws.Range("J38").Value = "A"
$ws.Range("K38").Value = 2.25
$ws.Range("L38").Value = 3.1
$ws.Range("M38").Value = 2.875
$ws.Range("N38").Value = 1.95
$ws.Range("O38").Value = 3
$ws.Range("P38").Value = 3.75
$ws.Range("Q38").Value = -0.5
$ws.Range("R38").Value = 2
$ws.Range("S38").Value = 1.8
$ws.Range("T38").Value = 2
$ws.Range("U38").Value = 1.85
$ws.Range("V38").Value = 1.95
$ws.Range("X38").Value = -1
$ws.Range("Y38").Value = 2.75
$ws.Range("AA38").Value = 0.8
$ws.Range("AB38").Value = 0.8500000000000001
$ws.Range("AC38").Value = -1

# Row 39
$ws.Range("B39").Value = 6979435
$ws.Range("G39").Value = "FK Radnik Surdulica"
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = "D"
$ws.Range("K39").Value = 1.4
$ws.Range("L39").Value = 3.75
$ws.Range("M39").Value = 7.5
$ws.Range("N39").Value = 1.5
$ws.Range("O39").Value = 3.6
$ws.Range("P39").Value = 6
$ws.Range("Q39").Value = -1
$ws.Range("R39").Value = 1.85
$ws.Range("S39").Value = 1.95
$ws.Range("T39").Value = 2.5
$ws.Range("U39").Value = 2
$ws.Range("V39").Value = 1.8
$ws.Range("X39").Value = 2.6
$ws.Range("Y39").Value = -1
$ws.Range("AA39").Value = 0.95
$ws.Range("AB39").Value = -1
$ws.Range("AC39").Value = 0.8

# Row 43
$ws.Range("B43").Value = 6978740
$ws.Range("F43").Value = "Red Star Belgrade"
$ws.Range("G43").Value = "FK Novi Pazar"
$ws.Range("H43").Value = 2
$ws.Range("I43").Value = 1
$ws.Range("K43").Value = 1.062
$ws.Range("L43").Value = 13
$ws.Range("M43").Value = 23
$ws.Range("N43").Value = 1.025
$ws.Range("O43").Value = 19
$ws.Range("P43").Value = 41
$ws.Range("Q43").Value = -3.75
$ws.Range("R43").Value = 1.825
$ws.Range("S43").Value = 1.975
$ws.Range("T43").Value = 4.5
$ws.Range("U43").Value = 1.975
$ws.Range("V43").Value = 1.825
$ws.Range("W43").Value = 0.02499999999999991
$ws.Range("Z43").Value = -1
$ws.Range("AA43").Value = 0.9750000000000001
$ws.Range("AC43").Value = 0.825

# Row 44
$ws.Range("B44").Value = 6979440
$ws.Range("F44").Value = "Javor Ivanjica"
$ws.Range("G44").Value = "Radnicki Nis"
$ws.Range("H44").Value = 1
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 2.3
$ws.Range("L44").Value = 3.2
$ws.Range("M44").Value = 2.875
$ws.Range("N44").Value = 2.5
$ws.Range("O44").Value = 3.25
$ws.Range("P44").Value = 2.6
$ws.Range("Q44").Value = 0
$ws.Range("R44").Value = 1.85
$ws.Range("S44").Value = 1.95
$ws.Range("T44").Value = 2.25
$ws.Range("U44").Value = 1.9
$ws.Range("V44").Value = 1.9
$ws.Range("W44").Value = 1.5
$ws.Range("Z44").Value = 0.8500000000000001
$ws.Range("AA44").Value = -1
$ws.Range("AC44").Value = 0.8999999999999999

# Row 46
$ws.Range("G46").Value = "FK Cukaricki"

# Row 48
$ws.Range("G48").Value = "Spartak Subotica"

# Row 49
$ws.Range("F49").Value = "Spartak Subotica"

# Row 50
$ws.Range("B50").Value = 6979449
$ws.Range("F50").Value = "FK Backa Topola"
$ws.Range("G50").Value = "FK Zeleznicar Pancevo"
$ws.Range("H50").Value = 6
$ws.Range("I50").Value = 3
$ws.Range("K50").Value = 1.25
$ws.Range("L50").Value = 5
$ws.Range("M50").Value = 9
$ws.Range("N50").Value = 1.285
$ws.Range("O50").Value = 4.5
$ws.Range("P50").Value = 9.5
$ws.Range("Q50").Value = -1.75
$ws.Range("T50").Value = 2.75
$ws.Range("U50").Value = 1.875
$ws.Range("V50").Value = 1.925
$ws.Range("W50").Value = 0.2849999999999999
$ws.Range("AB50").Value = 0.875

# Row 51
$ws.Range("B51").Value = 6979447
$ws.Range("F51").Value = "FK Vozdovac"
$ws.Range("G51").Value = "Mladost Lucani"
$ws.Range("H51").Value = 3
$ws.Range("I51").Value = 1
$ws.Range("K51").Value = 1.909
$ws.Range("L51").Value = 3.25
$ws.Range("M51").Value = 3.5
$ws.Range("N51").Value = 1.95
$ws.Range("O51").Value = 3.1
$ws.Range("P51").Value = 3.6
$ws.Range("Q51").Value = -0.5
$ws.Range("T51").Value = 2.25
$ws.Range("U51").Value = 1.9
$ws.Range("V51").Value = 1.9
$ws.Range("W51").Value = 0.95
$ws.Range("AB51").Value = 0.8999999999999999

# Row 52
$ws.Range("F52").Value = "FK Cukaricki"

# Row 59
$ws.Range("F59").Value = "FK Cukaricki"

# Row 64
$ws.Range("F64").Value = "FK Cukaricki"
$ws.Range("G64").Value = "Spartak Subotica"

# Row 77
$ws.Range("F77").Value = "Spartak Subotica"

# Row 78
$ws.Range("G78").Value = "FK Cukaricki"

# Row 79
$ws.Range("G79").Value = "Spartak Subotica"

# Row 81
$ws.Range("F81").Value = "FK Cukaricki"

# Row 84
$ws.Range("B84").Value = 6979481
$ws.Range("F84").Value = "Vojvodina"
$ws.Range("G84").Value = "FK Radnik Surdulica"
$ws.Range("H84").Value = 3
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = "H"
$ws.Range("K84").Value = 2.25
$ws.Range("L84").Value = 3
$ws.Range("M84").Value = 3
$ws.Range("N84").Value = 1.363
$ws.Range("O84").Value = 3.8
$ws.Range("P84").Value = 9
$ws.Range("Q84").Value = -1.25
$ws.Range("R84").Value = 1.825
$ws.Range("S84").Value = 1.975
$ws.Range("T84").Value = 2.5
$ws.Range("U84").Value = 1.9
$ws.Range("V84").Value = 1.9
$ws.Range("W84").Value = 0.363
$ws.Range("Y84").Value = -1
$ws.Range("Z84").Value = 0.825
$ws.Range("AA84").Value = -1
$ws.Range("AB84").Value = 0.8999999999999999

# Row 85
$ws.Range("B85").Value = 6979484
$ws.Range("F85").Value = "Mladost Lucani"
$ws.Range("G85").Value = "Radnicki Nis"
$ws.Range("H85").Value = 1
$ws.Range("I85").Value = 2
$ws.Range("J85").Value = "A"
$ws.Range("K85").Value = 2
$ws.Range("L85").Value = 3.25
$ws.Range("M85").Value = 3.25
$ws.Range("N85").Value = 2.55
$ws.Range("O85").Value = 3.2
$ws.Range("P85").Value = 2.55
$ws.Range("Q85").Value = 0
$ws.Range("R85").Value = 1.9
$ws.Range("S85").Value = 1.9
$ws.Range("T85").Value = 2.25
$ws.Range("U85").Value = 1.875
$ws.Range("V85").Value = 1.925
$ws.Range("W85").Value = -1
$ws.Range("Y85").Value = 1.55
$ws.Range("Z85").Value = -1
$ws.Range("AA85").Value = 0.8999999999999999
$ws.Range("AB85").Value = 0.875

# Row 90
$ws.Range("B90").Value = 6978747
$ws.Range("F90").Value = "IMT Novi Belgrade"
$ws.Range("G90").Value = "Red Star Belgrade"
$ws.Range("I90").Value = 2
$ws.Range("J90").Value = "A"
$ws.Range("K90").Value = 8
$ws.Range("L90").Value = 5.25
$ws.Range("M90").Value = 1.285
$ws.Range("N90").Value = 15
$ws.Range("O90").Value = 7.5
$ws.Range("P90").Value = 1.125
$ws.Range("Q90").Value = 2.25
$ws.Range("R90").Value = 1.975
$ws.Range("S90").Value = 1.825
$ws.Range("T90").Value = 3.5
$ws.Range("U90").Value = 1.825
$ws.Range("V90").Value = 1.975
$ws.Range("X90").Value = -1
$ws.Range("Y90").Value = 0.125
$ws.Range("Z90").Value = 0.9750000000000001
$ws.Range("AA90").Value = -1
$ws.Range("AC90").Value = 0.9750000000000001

# Row 91
$ws.Range("B91").Value = 6979491
$ws.Range("F91").Value = "Radnicki Nis"
$ws.Range("G91").Value = "FK Cukaricki"
$ws.Range("I91").Value = 1
$ws.Range("J91").Value = "D"
$ws.Range("K91").Value = 1.95
$ws.Range("L91").Value = 3.25
$ws.Range("M91").Value = 3.7
$ws.Range("N91").Value = 1.65
$ws.Range("O91").Value = 3.5
$ws.Range("P91").Value = 5
$ws.Range("Q91").Value = -0.75
$ws.Range("R91").Value = 1.825
$ws.Range("S91").Value = 1.975
$ws.Range("T91").Value = 2.5
$ws.Range("U91").Value = 2
$ws.Range("V91").Value = 1.8
$ws.Range("X91").Value = 2.5
$ws.Range("Y91").Value = -1
$ws.Range("Z91").Value = -1
$ws.Range("AA91").Value = 0.9750000000000001
$ws.Range("AC91").Value = 0.8

# Row 93
$ws.Range("F93").Value = "Spartak Subotica"

# Row 99
$ws.Range("F99").Value = "Spartak Subotica"

# Row 102
$ws.Range("F102").Value = "FK Cukaricki"

# Row 110
$ws.Range("G110").Value = "FK Cukaricki"

# Row 111
$ws.Range("G111").Value = "Spartak Subotica"

# Row 113
$ws.Range("G113").Value = "Spartak Subotica"

# Row 114
$ws.Range("G114").Value = "FK Cukaricki"

# Row 119
$ws.Range("G119").Value = "Spartak Subotica"

# Row 123
$ws.Range("F123").Value = "FK Cukaricki"

# Row 124
$ws.Range("B124").Value = 6979516
$ws.Range("F124").Value = "Partizan Belgrade"
$ws.Range("G124").Value = "Vojvodina"
$ws.Range("H124").Value = 3
$ws.Range("I124").Value = 1
$ws.Range("K124").Value = 1.5
$ws.Range("L124").Value = 4
$ws.Range("M124").Value = 5.5
$ws.Range("N124").Value = 1.444
$ws.Range("O124").Value = 4.2
$ws.Range("P124").Value = 6
$ws.Range("Q124").Value = -1.25
$ws.Range("R124").Value = 2.025
$ws.Range("S124").Value = 1.775
$ws.Range("T124").Value = 2.75
$ws.Range("W124").Value = 0.444
$ws.Range("Z124").Value = 1.025
$ws.Range("AB124").Value = 0.7749999999999999
$ws.Range("AC124").Value = -1

# Row 125
$ws.Range("B125").Value = 6979522
$ws.Range("F125").Value = "Mladost Lucani"
$ws.Range("G125").Value = "FK Zeleznicar Pancevo"
$ws.Range("H125").Value = 1
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 2.15
$ws.Range("L125").Value = 3.25
$ws.Range("M125").Value = 3.1
$ws.Range("N125").Value = 2.1
$ws.Range("O125").Value = 3.4
$ws.Range("P125").Value = 3.1
$ws.Range("Q125").Value = -0.25
$ws.Range("R125").Value = 1.85
$ws.Range("S125").Value = 1.95
$ws.Range("T125").Value = 2.25
$ws.Range("W125").Value = 1.1
$ws.Range("Z125").Value = 0.8500000000000001
$ws.Range("AB125").Value = -1
$ws.Range("AC125").Value = 1.025

# Row 127
$ws.Range("F127").Value = "Spartak Subotica"

# Row 132
$ws.Range("G132").Value = "FK Cukaricki"

# Row 134
$ws.Range("G134").Value = "Spartak Subotica"

# Row 142
$ws.Range("F142").Value = "FK Cukaricki"

# Row 145
$ws.Range("F145").Value = "Spartak Subotica"

# Row 148
$ws.Range("G148").Value = "FK Cukaricki"

# Row 151
$ws.Range("F151").Value = "Spartak Subotica"

# Row 154
$ws.Range("G154").Value = "Spartak Subotica"

# Row 161
$ws.Range("G161").Value = "FK Cukaricki"

# Row 164
$ws.Range("F164").Value = "Spartak Subotica"

# Row 166
$ws.Range("F166").Value = "FK Cukaricki"

# Row 167
$ws.Range("B167").Value = 6979547
$ws.Range("F167").Value = "FK Backa Topola"
$ws.Range("G167").Value = "FK Radnik Surdulica"
$ws.Range("K167").Value = 1.333
$ws.Range("L167").Value = 4.333
$ws.Range("M167").Value = 7.5
$ws.Range("N167").Value = 1.25
$ws.Range("O167").Value = 4.75
$ws.Range("P167").Value = 10
$ws.Range("Q167").Value = -1.5
$ws.Range("R167").Value = 1.85
$ws.Range("S167").Value = 1.95
$ws.Range("T167").Value = 2.5
$ws.Range("U167").Value = 1.8
$ws.Range("V167").Value = 2
$ws.Range("W167").Value = 0.25
$ws.Range("Z167").Value = -1
$ws.Range("AA167").Value = 0.95
$ws.Range("AC167").Value = 1

# Row 168
$ws.Range("B168").Value = 6979545
$ws.Range("F168").Value = "Radnicki Nis"
$ws.Range("G168").Value = "Javor Ivanjica"
$ws.Range("K168").Value = 2
$ws.Range("L168").Value = 3.25
$ws.Range("M168").Value = 3.25
$ws.Range("N168").Value = 1.727
$ws.Range("O168").Value = 3.3
$ws.Range("P168").Value = 4.2
$ws.Range("Q168").Value = -0.5
$ws.Range("R168").Value = 1.825
$ws.Range("S168").Value = 1.975
$ws.Range("T168").Value = 2.25
$ws.Range("U168").Value = 1.825
$ws.Range("V168").Value = 1.975
$ws.Range("W168").Value = 0.7270000000000001
$ws.Range("Z168").Value = 0.825
$ws.Range("AA168").Value = -1
$ws.Range("AC168").Value = 0.9750000000000001

# Row 172
$ws.Range("G172").Value = "Spartak Subotica"

# Row 176
$ws.Range("G176").Value = "FK Cukaricki"

# Row 179
$ws.Range("F179").Value = "Spartak Subotica"

# Row 185
$ws.Range("G185").Value = "FK Cukaricki"

# Row 188
$ws.Range("B188").Value = 7921658
$ws.Range("F188").Value = "FK Radnik Surdulica"
$ws.Range("G188").Value = "FK Radnicki 1923"
$ws.Range("K188").Value = 2.7
$ws.Range("M188").Value = 2.4
$ws.Range("N188").Value = 2.55
$ws.Range("O188").Value = 2.875
$ws.Range("P188").Value = 2.75
$ws.Range("R188").Value = 1.85
$ws.Range("S188").Value = 1.95
$ws.Range("U188").Value = 1.85
$ws.Range("V188").Value = 1.95
$ws.Range("Y188").Value = 1.75
$ws.Range("AA188").Value = 0.95
$ws.Range("AC188").Value = 0.95

# Row 189
$ws.Range("B189").Value = 6979565
$ws.Range("F189").Value = "FK Cukaricki"
$ws.Range("G189").Value = "FK Novi Pazar"
$ws.Range("K189").Value = 2.3
$ws.Range("M189").Value = 2.8
$ws.Range("N189").Value = 2.4
$ws.Range("O189").Value = 2.9
$ws.Range("P189").Value = 2.8
$ws.Range("R189").Value = 1.725
$ws.Range("S189").Value = 2.075
$ws.Range("U189").Value = 1.75
$ws.Range("V189").Value = 2.05
$ws.Range("Y189").Value = 1.8
$ws.Range("AA189").Value = 1.075
$ws.Range("AC189").Value = 1.05

# Row 190
$ws.Range("B190").Value = 6979566
$ws.Range("F190").Value = "Mladost Lucani"
$ws.Range("G190").Value = "Spartak Subotica"
$ws.Range("H190").Value = 1
$ws.Range("I190").Value = 0
$ws.Range("K190").Value = 3.4
$ws.Range("L190").Value = 3.3
$ws.Range("M190").Value = 1.95
$ws.Range("N190").Value = 2.7
$ws.Range("O190").Value = 3
$ws.Range("P190").Value = 2.45
$ws.Range("Q190").Value = 0
$ws.Range("R190").Value = 2
$ws.Range("S190").Value = 1.8
$ws.Range("T190").Value = 2
$ws.Range("U190").Value = 1.75
$ws.Range("V190").Value = 2.05
$ws.Range("W190").Value = 1.7
$ws.Range("Z190").Value = 1
$ws.Range("AA190").Value = -1
$ws.Range("AB190").Value = -1
$ws.Range("AC190").Value = 1.05

# Row 191
$ws.Range("B191").Value = 7921659
$ws.Range("F191").Value = "Vojvodina"
$ws.Range("G191").Value = "FK Vozdovac"
$ws.Range("H191").Value = 2
$ws.Range("I191").Value = 1
$ws.Range("K191").Value = 1.6
$ws.Range("L191").Value = 3.6
$ws.Range("M191").Value = 4.75
$ws.Range("N191").Value = 1.45
$ws.Range("O191").Value = 3.8
$ws.Range("P191").Value = 6
$ws.Range("Q191").Value = -1
$ws.Range("R191").Value = 1.8
$ws.Range("S191").Value = 2
$ws.Range("T191").Value = 2.5
$ws.Range("U191").Value = 1.85
$ws.Range("V191").Value = 1.95
$ws.Range("W191").Value = 0.45
$ws.Range("Z191").Value = 0
$ws.Range("AA191").Value = 0
$ws.Range("AB191").Value = 0.8500000000000001
$ws.Range("AC191").Value = -1

# Row 195
$ws.Range("B195").Value = 6979577
$ws.Range("F195").Value = "Spartak Subotica"
$ws.Range("G195").Value = "FK Cukaricki"
$ws.Range("K195").Value = 1.533
$ws.Range("L195").Value = 4
$ws.Range("M195").Value = 6
$ws.Range("N195").Value = 1.5
$ws.Range("O195").Value = 4
$ws.Range("P195").Value = 7
$ws.Range("Q195").Value = -1
$ws.Range("R195").Value = 1.8
$ws.Range("S195").Value = 2
$ws.Range("U195").Value = 1.8
$ws.Range("V195").Value = 2
$ws.Range("X195").Value = 3
$ws.Range("AA195").Value = 1
$ws.Range("AB195").Value = 0.8

# Row 196
$ws.Range("B196").Value = 6979578
$ws.Range("F196").Value = "FK Novi Pazar"
$ws.Range("G196").Value = "FK Vozdovac"
$ws.Range("K196").Value = 1.7
$ws.Range("L196").Value = 3.6
$ws.Range("M196").Value = 4.8
$ws.Range("N196").Value = 1.7
$ws.Range("O196").Value = 3.5
$ws.Range("P196").Value = 5
$ws.Range("Q196").Value = -0.75
$ws.Range("R196").Value = 1.9
$ws.Range("S196").Value = 1.9
$ws.Range("U196").Value = 1.775
$ws.Range("V196").Value = 2.025
$ws.Range("X196").Value = 2.5
$ws.Range("AA196").Value = 0.8999999999999999
$ws.Range("AB196").Value = 0.7749999999999999

# Row 203
$ws.Range("G203").Value = "Spartak Subotica"

# Row 204
$ws.Range("F204").Value = "FK Cukaricki"

# Row 210
$ws.Range("F210").Value = "Spartak Subotica"

# Row 217
$ws.Range("G217").Value = "FK Cukaricki"

# Row 219
$ws.Range("G219").Value = "Spartak Subotica"

# Row 223
$ws.Range("F223").Value = "FK Cukaricki"

# Row 227
$ws.Range("G227").Value = "Spartak Subotica"

# Row 228
$ws.Range("G228").Value = "FK Cukaricki"

# Row 232
$ws.Range("H232").Value = 1
$ws.Range("I232").Value = 1
$ws.Range("J232").Value = "D"
$ws.Range("N232").Value = 1.909
$ws.Range("O232").Value = 3.1
$ws.Range("P232").Value = 4.2
$ws.Range("R232").Value = 1.95
$ws.Range("S232").Value = 1.85
$ws.Range("W232").Value = -1
$ws.Range("X232").Value = 2.1
$ws.Range("Y232").Value = -1
$ws.Range("Z232").Value = -1
$ws.Range("AA232").Value = 0.8500000000000001
$ws.Range("AB232").Value = -0.5
$ws.Range("AC232").Value = 0.475

# Row 233
$ws.Range("H233").Value = 1
$ws.Range("I233").Value = 0
$ws.Range("J233").Value = "H"
$ws.Range("O233").Value = 2.9
$ws.Range("P233").Value = 3.5
$ws.Range("R233").Value = 1.85
$ws.Range("S233").Value = 1.95
$ws.Range("T233").Value = 2
$ws.Range("U233").Value = 1.875
$ws.Range("V233").Value = 1.925
$ws.Range("W233").Value = 1.15
$ws.Range("X233").Value = -1
$ws.Range("Y233").Value = -1
$ws.Range("Z233").Value = 0.8500000000000001
$ws.Range("AA233").Value = -1
$ws.Range("AB233").Value = -1
$ws.Range("AC233").Value = 0.925
